# Organisatorisches/Arbeitszeiten.xlsx
# Commit message: "added pdfCreator for logic separation"
#
# Adds two new time-tracking entries for Daniel (13 Nov 2018 and 14 Nov 2018)
# for the new task "Reworking Layers Separation", matching the formatting of
# the existing rows (row 55), and moves the sheet's active selection to the
# new last cell. The SUMIF-based totals (row 3) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date number format etc.) from the last existing row so
# the two new rows look the same as the rest of the table.
$ws.Range("A55").Copy()
$ws.Range("A56:A57").PasteSpecial(-4122)

# Row 56: 13 Nov 2018 (serial 43417), Daniel, Reworking Layers Separation, 4h
$ws.Range("A56").Value = 43417
$ws.Range("B56").Value = "Daniel"
$ws.Range("C56").Value = "Reworking Layers Separation"
$ws.Range("D56").Value = 4

# Row 57: 14 Nov 2018 (serial 43418), Daniel, Reworking Layers Separation, 3h
$ws.Range("A57").Value = 43418
$ws.Range("B57").Value = "Daniel"
$ws.Range("C57").Value = "Reworking Layers Separation"
$ws.Range("D57").Value = 3

# Move the active selection to follow the newly added data, matching the
# author's saved view state.
$ws.Range("D57").Select() | Out-Null
